$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values such as
# "0.9985" or "240.67" keep their original string representation
# instead of being auto-converted into floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.045.70"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "1.829.53"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("D4").Value = "0.9985"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "240.67"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").Value = "0.6215"
$ws.Range("E6").Value = "  -6.51%  "
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "44.52"
$ws.Range("E8").Value = "  +6.22%  "
$ws.Range("D9").Value = "0.07502"
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("D10").Value = "0.2914"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("D11").Value = "22.71"
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").Value = "0.07622"
$ws.Range("E12").Value = "  -1.73%  "
$ws.Range("D13").Value = "1.828.51"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").Value = "4.952"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").Value = "0.6633"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").Value = "82.15"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").Value = "0.000009124"
$ws.Range("E17").Value = "  +9.16%  "
$ws.Range("D18").Value = "5.993"
$ws.Range("E18").Value = "  -2.03%  "
$ws.Range("D19").Value = "29.039.30"
$ws.Range("E19").Value = "  -0.54%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "224.54"
$ws.Range("E20").Value = "  -1.59%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "12.32"
$ws.Range("E21").Value = "  -1.18%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "7.186"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "0.9998"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "159.48"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "8.409"
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "0.1358"
$ws.Range("E27").Value = "  -3.60%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "17.82"
$ws.Range("E28").Value = "  -0.94%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "1.498"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "4.034"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "4.049"
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "1.204"
$ws.Range("E32").Value = "  +1.03%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "0.05215"
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "1.830"
$ws.Range("E34").Value = "  -1.84%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.152"
$ws.Range("E35").Value = "  +1.24%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.7321"
$ws.Range("E36").Value = "  -1.86%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.642"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.270.52"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.748"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.01781"
$ws.Range("E40").Value = "  -1.04%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "6.320"
$ws.Range("E41").Value = "  +7.55%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "0.8946"
$ws.Range("E42").Value = "  -4.49%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "101.88"
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.974.49"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "0.5117"
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "63.35"
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.00000000119"
$ws.Range("E48").Value = "  -1.10%  "
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").Value = "0.3961"
$ws.Range("E49").Value = "  -1.28%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "1.671"
$ws.Range("E50").Value = "  -5.01%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "8.869"
$ws.Range("E51").Value = "  +0.54%  "
